# Auto-generated Excel COM-interop script to update cryptos.xlsx price/volume data
# per the commit 'Updated cryptos list on Sun Mar 19 08:39:28 UTC 2023 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "27.148.49"
$ws.Cells.Item(2, 5).Value = "  -0.81%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.783.82"
$ws.Cells.Item(3, 5).Value = "  -1.62%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  -0.27%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'336.89"
$ws.Cells.Item(5, 5).Value = "  -1.70%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  -0.19%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.3818"
$ws.Cells.Item(7, 5).Value = "  +0.15%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -2.39%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'48.09"
$ws.Cells.Item(9, 5).Value = "  -1.48%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07445"
$ws.Cells.Item(11, 5).Value = "  -3.62%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'1.001"
$ws.Cells.Item(12, 5).Value = "  -0.43%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "21.75"
$ws.Cells.Item(13, 5).Value = "  -1.30%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'6.433"

# Row 15
$ws.Cells.Item(15, 4).Value = "1.783.31"
$ws.Cells.Item(15, 5).Value = "  -1.12%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'7.092"
$ws.Cells.Item(16, 5).Value = "  -2.17%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'0.00001092"
$ws.Cells.Item(17, 5).Value = "  -2.36%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.06636"
$ws.Cells.Item(18, 5).Value = "  -1.21%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'83.43"
$ws.Cells.Item(19, 5).Value = "  -2.91%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.23%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'6.521"
$ws.Cells.Item(21, 5).Value = "  -0.46%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'17.42"
$ws.Cells.Item(22, 5).Value = "  -1.13%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "27.146.97"
$ws.Cells.Item(23, 5).Value = "  -0.89%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "12.25"
$ws.Cells.Item(24, 5).Value = "  -8.51%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.373"
$ws.Cells.Item(25, 5).Value = "  -3.82%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).Value = "'21.12"
$ws.Cells.Item(26, 5).Value = "  -4.05%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "LidoDAOToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(27, 4).Value = "'2.494"
$ws.Cells.Item(27, 5).Value = "  -6.63%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "1.452"
$ws.Cells.Item(28, 5).Value = "  -1.84%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'155.61"
$ws.Cells.Item(29, 5).Value = "  +1.37%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "1.980.97"
$ws.Cells.Item(30, 5).Value = "  -1.30%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "133.87"
$ws.Cells.Item(31, 5).Value = "  -1.70%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.985"
$ws.Cells.Item(32, 5).Value = "  -1.45%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'6.044"
$ws.Cells.Item(33, 5).Value = "  -4.39%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'0.08661"
$ws.Cells.Item(34, 5).Value = "  -0.90%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "13.07"
$ws.Cells.Item(35, 5).Value = "  -5.71%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "1.623"
$ws.Cells.Item(36, 5).Value = "  -4.96%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.6834"
$ws.Cells.Item(37, 5).Value = "  -2.02%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'5.382"
$ws.Cells.Item(38, 5).Value = "  -4.05%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.06293"
$ws.Cells.Item(39, 5).Value = "  -3.12%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -4.40%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.02317"
$ws.Cells.Item(41, 5).Value = "  -4.07%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.236"
$ws.Cells.Item(42, 5).Value = "  -4.68%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'8.379"
$ws.Cells.Item(43, 5).Value = "  -6.34%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'14.30"
$ws.Cells.Item(44, 5).Value = "  -2.58%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.6432"
$ws.Cells.Item(45, 5).Value = "  -1.22%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.9997"
$ws.Cells.Item(46, 5).Value = "  -0.24%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "3.855"
$ws.Cells.Item(47, 5).Value = "  -4.10%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'2.130"
$ws.Cells.Item(48, 5).Value = "  -2.18%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'131.46"
$ws.Cells.Item(49, 5).Value = "  -0.76%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.07101"
$ws.Cells.Item(50, 5).Value = "  -3.34%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'78.62"
$ws.Cells.Item(51, 5).Value = "  -2.33%  "

Write-Output "Updated cryptos list values"